# Update the "Metadata" sheet (sheet1) of the ValueSet workbook to the new
# publication state: version bump, status -> draft, new date, updated
# contact info, and a new "Jurisdiction" row inserted before "Description".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- simple value updates -------------------------------------------------
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- insert a new "Jurisdiction" row before "Description" ----------------
# Read the existing rows 12-15 (Description, Purpose, Copyright, Immutable)
# using Value2 (Value has a read quirk in this host) before they get
# overwritten, then shift them down one row to make room for the new row.
$a12 = $ws.Range("A12").Value2
$b12 = $ws.Range("B12").Value2
$a13 = $ws.Range("A13").Value2
$b13 = $ws.Range("B13").Value2
$a14 = $ws.Range("A14").Value2
$b14 = $ws.Range("B14").Value2
$a15 = $ws.Range("A15").Value2
$b15 = $ws.Range("B15").Value2

# Copy row 15's formatting down onto the new row 16 before we touch values.
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

$ws.Range("A16").Value = $a15
$ws.Range("B16").Value = $b15
$ws.Range("A15").Value = $a14
$ws.Range("B15").Value = $b14
$ws.Range("A14").Value = $a13
$ws.Range("B14").Value = $b13
$ws.Range("A13").Value = $a12
$ws.Range("B13").Value = $b12

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

$excel.CutCopyMode = 0
